# Auto-generated PowerShell Excel COM-interop script
# Applies numeric value updates (and a few cell clears) to the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets per the target diff.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 12044.667
$ws.Range("I34").Value = 4653.6
$ws.Range("J34").Value = 49000
$ws.Range("K34").Value = 4653.6
$ws.Range("L34").Value = 49000
$ws.Range("M34").Value = -4450.6
$ws.Range("N34").Value = -49406
$ws.Range("H36").Value = 12044.667
$ws.Range("I36").Value = 4653.6
$ws.Range("J36").Value = 49000
$ws.Range("K36").Value = 4653.6
$ws.Range("L36").Value = 49000
$ws.Range("M36").Value = -3938.6
$ws.Range("N36").Value = -50430
$ws.Range("H80").Value = 3574.2727
$ws.Range("J80").Value = 3739.6
$ws.Range("L80").Value = 11218.8
$ws.Range("N80").Value = -13214.8
$ws.Range("H83").Value = 3574.2727
$ws.Range("J83").Value = 3739.6
$ws.Range("L83").Value = 33656.4
$ws.Range("N83").Value = -43640.4
$ws.Range("H127").Value = 1569
$ws.Range("I127").Value = 990
$ws.Range("J127").Value = 1858.5
$ws.Range("K127").Value = 2970
$ws.Range("L127").Value = 5575.5
$ws.Range("M127").Value = 1990
$ws.Range("N127").Value = -15495.5
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 59999.5
$ws.Range("J136").Value = 59999.5
$ws.Range("L136").Value = 59999.5
$ws.Range("N136").Value = -70199.5
$ws.Range("H137").Value = 1078
$ws.Range("I137").Value = 1078
$ws.Range("K137").Value = 3234
$ws.Range("M137").Value = -684
$ws.Range("H138").Value = 2291.4375
$ws.Range("I138").Value = 1264.6666
$ws.Range("J138").Value = 2528.3845
$ws.Range("K138").Value = 3793.9998
$ws.Range("L138").Value = 7585.1535
$ws.Range("M138").Value = 1346.0002
$ws.Range("N138").Value = -17865.1535
$ws.Range("N135").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2675.6924
$ws.Range("I45").Value = 2806.0833
$ws.Range("J45").Value = 1111
$ws.Range("K45").Value = 2806.0833
$ws.Range("L45").Value = 1111
$ws.Range("M45").Value = -2429.0833
$ws.Range("N45").Value = -1865
$ws.Range("H61").Value = 1849.1666
$ws.Range("I61").Value = 1849.1666
$ws.Range("K61").Value = 1849.1666
$ws.Range("M61").Value = -1637.1666
$ws.Range("H74").Value = 9501.583000000001
$ws.Range("I74").Value = 9456.362999999999
$ws.Range("K74").Value = 9456.362999999999
$ws.Range("M74").Value = -8582.362999999999
$ws.Range("H77").Value = 9501.583000000001
$ws.Range("I77").Value = 9456.362999999999
$ws.Range("K77").Value = 47281.815
$ws.Range("M77").Value = -42913.815
$ws.Range("H110").Value = 3653.3572
$ws.Range("I110").Value = 2186.182
$ws.Range("K110").Value = 2186.182
$ws.Range("M110").Value = -141.1819999999998
$ws.Range("H136").Value = 1849.1666
$ws.Range("I136").Value = 1849.1666
$ws.Range("K136").Value = 5547.4998
$ws.Range("M136").Value = -2997.4998

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5651.3076
$ws.Range("I94").Value = 5497.778
$ws.Range("J94").Value = 5996.75
$ws.Range("K94").Value = 5497.778
$ws.Range("L94").Value = 5996.75
$ws.Range("M94").Value = -5046.778
$ws.Range("N94").Value = -6898.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1839.4445
$ws.Range("I31").Value = 1501.5
$ws.Range("K31").Value = 1501.5
$ws.Range("M31").Value = -1206.5
$ws.Range("H34").Value = 1839.4445
$ws.Range("I34").Value = 1501.5
$ws.Range("K34").Value = 1501.5
$ws.Range("M34").Value = -1299.5
$ws.Range("H134").Value = 9885.799999999999
$ws.Range("I134").Value = 10133
$ws.Range("K134").Value = 30399
$ws.Range("M134").Value = -27864

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2831
$ws.Range("H23").Value = 1007
$ws.Range("I23").Value = 1009.5
$ws.Range("K23").Value = 3028.5
$ws.Range("M23").Value = -2793.5
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -2898
$ws.Range("H129").Value = 2669.75
$ws.Range("I129").Value = 1399
$ws.Range("K129").Value = 4197
$ws.Range("M129").Value = 803
$ws.Range("M4").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("N27").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 35986.25
$ws.Range("J141").Value = 35986.25
$ws.Range("L141").Value = 35986.25
$ws.Range("N141").Value = -46346.25

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4641.5386
$ws.Range("I46").Value = 1445
$ws.Range("K46").Value = 1445
$ws.Range("M46").Value = -1257
$ws.Range("H132").Value = 5436.5
$ws.Range("I132").Value = 4500
$ws.Range("K132").Value = 13500
$ws.Range("M132").Value = -10970
$ws.Range("H136").Value = 3579.5
$ws.Range("I136").Value = 3579.5
$ws.Range("K136").Value = 10738.5
$ws.Range("M136").Value = -8188.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2034
$ws.Range("I81").Value = 1913.25
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 3826.5
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -2765.5
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 2034
$ws.Range("I84").Value = 1913.25
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 19132.5
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -13828.5
$ws.Range("N84").Value = -40608
$ws.Range("H132").Value = 2604.2
$ws.Range("I132").Value = 1507.6666
$ws.Range("K132").Value = 4522.9998
$ws.Range("M132").Value = -1992.9998
$ws.Range("H136").Value = 5564.4
$ws.Range("I136").Value = 4830.625
$ws.Range("J136").Value = 8499.5
$ws.Range("K136").Value = 14491.875
$ws.Range("L136").Value = 25498.5
$ws.Range("M136").Value = -11941.875
$ws.Range("N136").Value = -30598.5
$ws.Range("H140").Value = 29998.334
$ws.Range("J140").Value = 29998.334
$ws.Range("L140").Value = 29998.334
$ws.Range("N140").Value = -40358.334
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
